$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.241.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.025.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.023.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.21%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.260.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.527.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.022.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.06%  "
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000101"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("E33").Value = "  +5.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.990"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.45%  "
$ws.Range("E39").Value = "  -6.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.312"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "385.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.735.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.28%  "
